# edit.ps1 - applies the "Removed register option" commit to the Architecture Notebook template.
#
# Summary of changes performed:
#  1. Remove the stray _GoBack bookmark that sat right after the title "ACAD Oracle".
#  2. Extend the "Persistency" bullet with a new trailing sentence describing how
#     persistence is handled (SQL Server). Word leaves its autosave "last edit"
#     _GoBack bookmark wrapping the freshly typed words, so we recreate it there.
#  3. A handful of other bullet paragraphs got small grammar/spell-check passes
#     re-run over them (their runs collapse into a single run and the stray
#     <w:proofErr> tags disappear) - we reproduce that by doing a self
#     find&replace over each affected sentence, which is exactly what happens
#     when Word's proofer re-evaluates a paragraph after it is touched.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: not found -> $old"
    }
}

# --- 1. Drop the original _GoBack bookmark (right after the document title). ---
$d.Bookmarks("_GoBack").Delete()

# --- 2. Collapse the runs / drop proofErr spans by replacing each sentence with itself. ---
Replace-Text "The development team shall understand the concept of MVC, or shall have experience working in a multi-layered architecture" "The development team shall understand the concept of MVC, or shall have experience working in a multi-layered architecture"

Replace-Text "A SQL Server Express database shall be installed in the database server." "A SQL Server Express database shall be installed in the database server."

Replace-Text "The system shall be connected to the network in order to be used by the client. If either the database server and/or application server cannot have network, the system will not be usable." "The system shall be connected to the network in order to be used by the client. If either the database server and/or application server cannot have network, the system will not be usable."

Replace-Text "Web application shall be implemented using the Microsoft MVC framework 4.0 with C# language. The framework is robust and reliable, and also have an easy to implement authentication functionality and session management." "Web application shall be implemented using the Microsoft MVC framework 4.0 with C# language. The framework is robust and reliable, and also have an easy to implement authentication functionality and session management."

Replace-Text "The data layer shall be implemented using the Microsoft Entity Framework. It enables quick setup and easy implementation to communicate with the database." "The data layer shall be implemented using the Microsoft Entity Framework. It enables quick setup and easy implementation to communicate with the database."

Replace-Text "The Controllers in the MVC project can reference the repositories. But only their interfaces. The implementation shall be injected using the injection pattern. The " "The Controllers in the MVC project can reference the repositories. But only their interfaces. The implementation shall be injected using the injection pattern. The "

Replace-Text " library shall be used to inject the dependencies in the controller" " library shall be used to inject the dependencies in the controller"

Replace-Text "This way, in the controller" "This way, in the controller"

Replace-Text " object and the implementation will be injected by simple injector. " " object and the implementation will be injected by simple injector. "

Replace-Text "can be accessed through a browser and doesn" "can be accessed through a browser and doesn"

Replace-Text "to display useful user friendly messages." "to display useful user friendly messages."

Replace-Text " Log4net shall be used for logging. " " Log4net shall be used for logging. "

Replace-Text " will be used for it. The documentation for " " will be used for it. The documentation for "

Replace-Text " They will not be persisted, they will be evaluated in runtime and provided by the user when generating the suggested disciplines." " They will not be persisted, they will be evaluated in runtime and provided by the user when generating the suggested disciplines."

Replace-Text "Classes: will define the day/time the discipline can be coursed." "Classes: will define the day/time the discipline can be coursed."

# --- 3. Extend the "Persistency" sentence with the new text about the SQL Server backing store,
#        matching the new sentence trailing where Word's _GoBack bookmark ends up. ---
Replace-Text "to save the data used to suggest the disciplines to be coursed." "to save the data used to suggest the disciplines to be coursed. This is throught the Sql Server relational database."

